$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLAN_INFO")

# Simplify the route: drop the extra stops (columns G:K), and move the last
# stop's value into column F so the route becomes a short round trip.
$ws.Range("F2").Value = "HKNW"
$ws.Range("G2:K2").ClearContents()

$ws.Range("F3:J3").ClearContents()
$ws.Range("F4:J4").ClearContents()
$ws.Range("F5:J5").ClearContents()
$ws.Range("F6:J6").ClearContents()

$ws.Range("B7").Value = 100
$ws.Range("E7").Value = 110
$ws.Range("F7:J7").ClearContents()

# Crew now uses short codes instead of full names.
$ws.Range("B13").Value = "AS"
$ws.Range("C13").Value = "BM"
$ws.Range("A13").Value = "CREW (PIC followed by SIC)"

# TOF changed.
$ws.Range("B14").Value = 1200

$ws.Range("F7").Select()
